# fix latency units in report sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to the Read Latency min/max/average columns (I, J, K) for data rows 3-15
for ($row = 3; $row -le 15; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $orig = $cell.Text
        $cell.Value = "$orig msec"
    }
}
